$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A22").Value = "Hi there! I’m your sales agent from Creer Infotech. I’ve reached out to share some exciting offers on our latest products. Can I take a few minutes to tell you about them?"
$ws.Range("B22").Value = "earbuds"
$ws.Range("C22").Value = "neutral"
$ws.Range("D22").Value = "I'm sorry, I didn’t catch that."
$ws.Range("E22").Value = "2025-10-31 18:56:43"

$ws.Range("A23").Value = "I'm sorry, I didn’t catch that."
$ws.Range("B23").Value = "Bluetooth earbuds"
$ws.Range("C23").Value = "neutral"
$ws.Range("D23").Value = "Great choice! I’ve sent the link of Bluetooth Earbuds to your phone number ending with 0234. Thank you for your time! I really appreciate it. If you need anything, feel free to contact us.Our Contact number is 20251."
$ws.Range("E23").Value = "2025-10-31 18:57:09"

$ws.Range("A24").Value = "Hi there! I’m your sales agent from Creer Infotech. I’ve reached out to share some exciting offers on our latest products. Can I take a few minutes to tell you about them?"
$ws.Range("B24").Value = "Bluetooth earbuds"
$ws.Range("C24").Value = "neutral"
$ws.Range("D24").Value = "Great choice! I’ve sent the link of Bluetooth Earbuds to your phone number ending with 0234. Thank you for your time! I really appreciate it. If you need anything, feel free to contact us.Our Contact number is 20251."
$ws.Range("E24").Value = "2025-10-31 19:00:36"

$ws.Range("A25").Value = "Hi there! I’m your sales agent from Creer Infotech. I’ve reached out to share some exciting offers on our latest products. Can I take a few minutes to tell you about them?"
$ws.Range("B25").Value = "Bluetooth earbuds"
$ws.Range("C25").Value = "neutral"
$ws.Range("D25").Value = "Great choice! I’ve sent the link of Bluetooth Earbuds to your phone number ending with 0234. Thank you for your time! I really appreciate it. If you need anything, feel free to contact us. Our Contact number is 20251."
$ws.Range("E25").Value = "2025-10-31 19:13:09"

$ws.Range("A26").Value = "Hi there! I’m your sales agent from Creer Infotech. I’ve reached out to share some exciting offers on our latest products. Can I take a few minutes to tell you about them?"
$ws.Range("B26").Value = "Bluetooth earbuds"
$ws.Range("C26").Value = "neutral"
$ws.Range("D26").Value = "Great choice! I’ve sent the link of Bluetooth Earbuds to your phone number ending with 0234. Thank you for your time! I really appreciate it. If you need anything, feel free to contact us. Our Contact number is 20251."
$ws.Range("E26").Value = "2025-10-31 19:22:07"

$ws.Range("A27").Value = "Hi there! I’m your sales agent from Creer Infotech. I’ve reached out to share some exciting offers on our latest products. Can I take a few minutes to tell you about them?"
$ws.Range("B27").Value = "Bluetooth earbuds"
$ws.Range("C27").Value = "neutral"
$ws.Range("D27").Value = "Great choice! I’ve sent the link of Bluetooth Earbuds to your phone number ending with 0234. Thank you for your time! I really appreciate it. If you need anything, feel free to contact us. Our Contact number is 20251."
$ws.Range("E27").Value = "2025-10-31 19:34:15"

$ws.Range("A28").Value = "Hi there! I’m your sales agent from Creer Infotech. I’ve reached out to share some exciting offers on our latest products. Can I take a few minutes to tell you about them?"
$ws.Range("B28").Value = "Bluetooth earbuds"
$ws.Range("C28").Value = "neutral"
$ws.Range("D28").Value = "Great choice! I’ve sent the link of Bluetooth Earbuds to your phone number ending with 0234. Thank you for your time! I really appreciate it. If you need anything, feel free to contact us. Our Contact number is 20251."
$ws.Range("E28").Value = "2025-10-31 19:39:39"
